$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying TPM computation was rerun and the "Inflammatory-Mac"
# target-cluster row (original sheet row 4) no longer appears in the
# result set. Delete that row outright; the two rows below it (MuSCs and
# Resolving-Mac targets) shift up into rows 4 and 5, and every recomputed
# metric in columns E:T gets its new value.
$ws.Rows(4).Delete() | Out-Null

function Set-RowValues($rowIndex, $values) {
    $col = 5  # column E
    foreach ($v in $values) {
        $ws.Cells.Item($rowIndex, $col).Value = $v
        $col++
    }
}

# Row 2: MuSCs -> Edn3/Ednra -> ECs
Set-RowValues 2 @(2, 1, 0.307034, 0.6140680000000001, 1, 1, 2, 0.6666666666666666, 0.7195943333333332, 2.158783, 0.03204779321415739, 0.03468223907394029, 0.2209399265406667, 1.325639559244, 0.03204779321415739, 0.03468223907394029)

# Row 3: MuSCs -> Edn3/Ednra -> FAPs
Set-RowValues 3 @(2, 1, 0.307034, 0.6140680000000001, 1, 1, 3, 1, 16.53477466666667, 49.60432400000001, 0.7363913455312854, 0.7969254084681946, 5.076738005005334, 30.46042803003201, 0.7363913455312854, 0.7969254084681946)

# Row 4: MuSCs -> Edn3/Ednra -> MuSCs (this row used to be row 5 before the delete)
Set-RowValues 4 @(2, 1, 0.307034, 0.6140680000000001, 1, 1, 2, 1, 5.116736, 10.233472, 0.2278785277530463, 0.1644073176694804, 1.571011921024, 6.284047684096, 0.2278785277530463, 0.1644073176694804)

# Row 5: MuSCs -> Edn3/Ednra -> Resolving-Mac (this row used to be row 6 before the delete)
Set-RowValues 5 @(2, 1, 0.307034, 0.6140680000000001, 1, 1, 2, 0.6666666666666666, 0.08268233333333334, 0.248047, 0.003682333501510851, 0.003985034788384784, 0.02538628753266667, 0.152317725196, 0.003682333501510851, 0.003985034788384784)
